# Update the "remn_amt" time series on each of the 8 worksheets:
#   - row 105 (date 45968) had a placeholder B value of 0; fill in the real amount
#   - append 5 new rows (106-110) for dates 45971-45975 with their amounts
#
# Data per sheet (in workbook/tab order), B105 value + the 5 new (date, amount) rows.

$wb = $excel.ActiveWorkbook

$sheetData = @(
    @{ B105 = 10330162; Rows = @(
            @(45971, 10684685),
            @(45972, 11033960),
            @(45973, 10883200),
            @(45974, 10769447),
            @(45975, 10121585)
        ) },
    @{ B105 = 13831084; Rows = @(
            @(45971, 14607791),
            @(45972, 14548260),
            @(45973, 14410223),
            @(45974, 14120461),
            @(45975, 11744953)
        ) },
    @{ B105 = 3419631; Rows = @(
            @(45971, 3428991),
            @(45972, 3666127),
            @(45973, 3656561),
            @(45974, 3625031),
            @(45975, 3437693)
        ) },
    @{ B105 = 997504; Rows = @(
            @(45971, 1047257),
            @(45972, 1048409),
            @(45973, 1008852),
            @(45974, 1057005),
            @(45975, 991615)
        ) },
    @{ B105 = 1538375; Rows = @(
            @(45971, 1555025),
            @(45972, 1642945),
            @(45973, 1624142),
            @(45974, 1624112),
            @(45975, 1525563)
        ) },
    @{ B105 = 1762113; Rows = @(
            @(45971, 1744502),
            @(45972, 1883229),
            @(45973, 1904385),
            @(45974, 1940920),
            @(45975, 1800878)
        ) },
    @{ B105 = 284818; Rows = @(
            @(45971, 299301),
            @(45972, 305333),
            @(45973, 302312),
            @(45974, 309523),
            @(45975, 301426)
        ) },
    @{ B105 = 288357; Rows = @(
            @(45971, 321156),
            @(45972, 339672),
            @(45973, 341932),
            @(45974, 337354),
            @(45975, 321920)
        ) }
)

for ($i = 0; $i -lt $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i + 1)
    $data = $sheetData[$i]

    # Fill in the real amount that replaces the placeholder 0 in B105.
    $ws.Cells.Item(105, 2).Value = $data.B105

    # The date column (A) uses a custom date/time number format; grab it from
    # the last existing data row so the appended rows match formatting.
    $dateFormat = $ws.Range("A105").NumberFormat

    $r = 106
    foreach ($row in $data.Rows) {
        $ws.Range("A" + $r).NumberFormat = $dateFormat
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        $r++
    }
}
